$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.937071323394775
$ws.Range("B1").Value = 6.741384983062744
$ws.Range("C1").Value = 5.762762069702148
$ws.Range("D1").Value = 4.924643039703369
$ws.Range("E1").Value = 2.070333003997803
